$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A75").Value = "2025-04-29 11:17:30"
$ws.Range("B75").Value = 254
